$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1681159420289855
$ws.Range("C2").Value = 0.6086956521739131
$ws.Range("J2").Value = 0.008695652173913044
$ws.Range("P2").Value = 0.136231884057971
$ws.Range("S2").Value = 0.0782608695652174
$ws.Range("C3").Value = 0.02727272727272727
$ws.Range("J3").Value = 0.02727272727272727
$ws.Range("P3").Value = 0.7818181818181819
$ws.Range("S3").Value = 0.1636363636363636
$ws.Range("J4").Value = 0.03636363636363636
$ws.Range("P4").Value = 0.7818181818181819
$ws.Range("S4").Value = 0.1818181818181818
$ws.Range("B6").Value = 0.05976095617529881
$ws.Range("D6").Value = 0.01195219123505976
$ws.Range("F6").Value = 0.0796812749003984
$ws.Range("J6").Value = 0.2549800796812749
$ws.Range("O6").Value = 0.00398406374501992
$ws.Range("Q6").Value = 0.1195219123505976
$ws.Range("R6").Value = 0.0796812749003984
$ws.Range("S6").Value = 0.3904382470119522
$ws.Range("B7").Value = 0.09289617486338798
$ws.Range("D7").Value = 0.0273224043715847
$ws.Range("F7").Value = 0.07650273224043716
$ws.Range("J7").Value = 0.09836065573770492
$ws.Range("O7").Value = 0.00546448087431694
$ws.Range("Q7").Value = 0.1748633879781421
$ws.Range("R7").Value = 0.07103825136612021
$ws.Range("S7").Value = 0.453551912568306
$ws.Range("B8").Value = 0.1015625
$ws.Range("D8").Value = 0.01953125
$ws.Range("F8").Value = 0.0859375
$ws.Range("J8").Value = 0.109375
$ws.Range("O8").Value = 0.01953125
$ws.Range("Q8").Value = 0.14453125
$ws.Range("R8").Value = 0.10546875
$ws.Range("S8").Value = 0.4140625
$ws.Range("B9").Value = 0.1504854368932039
$ws.Range("D9").Value = 0.01941747572815534
$ws.Range("F9").Value = 0.07766990291262135
$ws.Range("J9").Value = 0.1067961165048544
$ws.Range("O9").Value = 0.009708737864077669
$ws.Range("Q9").Value = 0.1796116504854369
$ws.Range("R9").Value = 0.0825242718446602
$ws.Range("S9").Value = 0.3737864077669903
$ws.Range("B10").Value = 0.1258327165062916
$ws.Range("D10").Value = 0.02442635085122132
$ws.Range("E10").Value = 0.0007401924500370096
$ws.Range("F10").Value = 0.06439674315321983
$ws.Range("J10").Value = 0.1117690599555885
$ws.Range("O10").Value = 0.01332346410066617
$ws.Range("Q10").Value = 0.1998519615099926
$ws.Range("R10").Value = 0.09030347890451518
$ws.Range("S10").Value = 0.3693560325684678
$ws.Range("G11").Value = 0.1407942238267148
$ws.Range("J11").Value = 0.1046931407942238
$ws.Range("K11").Value = 0.1985559566787004
$ws.Range("L11").Value = 0.5379061371841155
$ws.Range("S11").Value = 0.01805054151624549
$ws.Range("G12").Value = 0.7712418300653595
$ws.Range("J12").Value = 0.1568627450980392
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.03267973856209151
$ws.Range("S12").Value = 0.03267973856209151
$ws.Range("G13").Value = 0.6862745098039216
$ws.Range("J13").Value = 0.2745098039215687
$ws.Range("S13").Value = 0.0392156862745098
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.0187793427230047
$ws.Range("H15").Value = 0.2300469483568075
$ws.Range("I15").Value = 0.07511737089201878
$ws.Range("J15").Value = 0.323943661971831
$ws.Range("K15").Value = 0.04225352112676056
$ws.Range("M15").Value = 0.01408450704225352
$ws.Range("O15").Value = 0.04694835680751173
$ws.Range("S15").Value = 0.2488262910798122
$ws.Range("F16").Value = 0.02371541501976284
$ws.Range("H16").Value = 0.1620553359683795
$ws.Range("I16").Value = 0.07114624505928854
$ws.Range("J16").Value = 0.4071146245059288
$ws.Range("K16").Value = 0.1225296442687747
$ws.Range("M16").Value = 0.02766798418972332
$ws.Range("N16").Value = 0.003952569169960474
$ws.Range("O16").Value = 0.04743083003952569
$ws.Range("S16").Value = 0.1343873517786561
$ws.Range("F17").Value = 0.01801801801801802
$ws.Range("H17").Value = 0.1621621621621622
$ws.Range("I17").Value = 0.09234234234234234
$ws.Range("J17").Value = 0.4572072072072072
$ws.Range("K17").Value = 0.06531531531531531
$ws.Range("M17").Value = 0.02027027027027027
$ws.Range("N17").Value = 0.002252252252252252
$ws.Range("O17").Value = 0.08108108108108109
$ws.Range("S17").Value = 0.1013513513513514
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.1785714285714286
$ws.Range("I18").Value = 0.08035714285714286
$ws.Range("J18").Value = 0.4151785714285715
$ws.Range("K18").Value = 0.08482142857142858
$ws.Range("M18").Value = 0.008928571428571428
$ws.Range("N18").Value = 0.004464285714285714
$ws.Range("O18").Value = 0.08035714285714286
$ws.Range("S18").Value = 0.1160714285714286
$ws.Range("F19").Value = 0.01409495548961424
$ws.Range("H19").Value = 0.2336795252225519
$ws.Range("I19").Value = 0.08234421364985163
$ws.Range("J19").Value = 0.3805637982195846
$ws.Range("K19").Value = 0.09495548961424333
$ws.Range("M19").Value = 0.02299703264094955
$ws.Range("N19").Value = 0.000741839762611276
$ws.Range("O19").Value = 0.05712166172106825
$ws.Range("S19").Value = 0.1135014836795252
